# Append a new data row (row 18) to the Adafruit IO export sheet, mirroring
# the existing rows (timestamp / feed key / value / lat / lon / elevation).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 18

# Column C holds a numeric-looking reading ("25"); force it to stay text
# (matching every other row in the sheet) instead of being auto-converted
# to a number.
$ws.Range("C$row").NumberFormat = "@"

$ws.Range("A$row").Value = "2024-09-25T18:06:40Z"
$ws.Range("B$row").Value = "temperature"
$ws.Range("C$row").Value = "25"
$ws.Range("D$row").Value = "N/A"
$ws.Range("E$row").Value = "N/A"
$ws.Range("F$row").Value = "N/A"
